# ---------------------------------------------------------------------------
# Servir.xlsx edit script
#
# Summary of the change (see commit message):
#  - The "lookup" sheet (tab name "Sheet1", file xl/worksheets/sheet2.xml)
#    gets 15 new PU/GRE/CG/#### codes (column A) together with their
#    matching UUID ids (column B), appended as rows 1119-1133. They are
#    written column-by-column (all codes, then all ids) to mirror how the
#    author pasted them in ("segmenté el origen ... por columnas").
#  - The "Hoja1" sheet (file xl/worksheets/sheet1.xml) is the main table:
#      * row 2/3/4 get new household codes (B) / dates (C) / delivery
#        location (E) values
#      * 3 new rows (5,6,7) are appended, following the same pattern
#        (household code / date / motive / delivery location / computed
#        code / VLOOKUP id / signature flag)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsLookup = $wb.Worksheets.Item("Sheet1")   # xl/worksheets/sheet2.xml
$wsHoja1  = $wb.Worksheets.Item("Hoja1")    # xl/worksheets/sheet1.xml

# ---------------------------------------------------------------------------
# 1) Lookup table ("Sheet1"): append the 15 new code/id pairs as rows
#    1119..1133, right after the existing last row (1118).
# ---------------------------------------------------------------------------

$newCodes = @(
    "PU/GRE/CG/0308/01",
    "PU/GRE/CG/0308/02",
    "PU/GRE/CG/0308/03",
    "PU/GRE/CG/0308/04",
    "PU/GRE/CG/0309/01",
    "PU/GRE/CG/0309/02",
    "PU/GRE/CG/0309/03",
    "PU/GRE/CG/0309/04",
    "PU/GRE/CG/0310/01",
    "PU/GRE/CG/0310/02",
    "PU/GRE/CG/0310/03",
    "PU/GRE/CG/0310/04",
    "PU/GRE/CG/0311/01",
    "PU/GRE/CG/0311/02",
    "PU/GRE/CG/0311/03"
)

$newIds = @(
    "64277d0e-97da-446d-a687-0d7ea0fe64b3",
    "8875a779-0a18-4aed-86ba-1b03b378e439",
    "e3a0cbf8-6de8-4627-9384-f02a496569d4",
    "33133711-ae3f-4166-959b-4b99ab79d07d",
    "c128fd59-88ae-4ebe-9569-15be7e13b096",
    "7e1d7969-7da4-4889-b4f6-aaeae9d46961",
    "19036049-ad12-456e-8fe0-9527f021ce80",
    "fbc2cb3c-477a-4933-95ba-09d4c1d2c4cb",
    "11e39041-f9fd-45b6-94b8-c301b9910979",
    "1769f7e0-f29c-4795-83d6-7a8838da5757",
    "b70bc750-687b-4d7a-80fc-2210aa945b09",
    "e7e0bf58-c6b2-42dc-a78f-0e9b8ae98e47",
    "fc8be43a-4d54-457b-a84b-406a829a0135",
    "381e38e5-745b-4249-93de-de7da7aae679",
    "e9d0910c-e2cf-4ba4-b694-b66edeaf45d3"
)

$startRow = 1119

for ($i = 0; $i -lt $newCodes.Length; $i++) {
    $wsLookup.Cells.Item($startRow + $i, 1).Value = $newCodes[$i]
}
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $wsLookup.Cells.Item($startRow + $i, 2).Value = $newIds[$i]
}

# ---------------------------------------------------------------------------
# 2) Main sheet ("Hoja1"): update existing rows 2-4, then append rows 5-7.
# ---------------------------------------------------------------------------

# --- Row 2: household code 111 -> 310 ---
$wsHoja1.Range("B2").Value = 310

# --- Row 3: household code 166 -> 311, date moves a few days, delivery
#            location switches from "En el SAI" to the new "En el Hogar" ---
$wsHoja1.Range("B3").Value = 311
$wsHoja1.Range("C3").Value = 45061
$wsHoja1.Range("E3").Value = "En el Hogar"

# --- Row 4: household code 217 -> 116 ---
$wsHoja1.Range("B4").Value = 116

# --- Rows 5-7: brand-new entries. Clone row 4's formatting first so the
#     B/C/H columns keep their existing number-format styles (PU/GRE/CG
#     code style, date style, integer style), then fill in the values and
#     formulas cell by cell (mirrors F2/G2's formulas exactly, with the
#     relative row references Excel would generate on fill-down). ---
$wsHoja1.Range("A4:H4").Copy()
$wsHoja1.Range("A5:H7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 5
$wsHoja1.Range("A5").Value = 3
$wsHoja1.Range("B5").Value = 290
$wsHoja1.Range("C5").Value = 45058
$wsHoja1.Range("D5").Value = "Seguimiento"
$wsHoja1.Range("E5").Value = "En el SAI"
$wsHoja1.Range("F5").Formula = '=CONCATENATE(TEXT(B5,"0000"),"/",TEXT(H5,"00"))'
$wsHoja1.Range("G5").Formula = '=VLOOKUP(CONCATENATE("PU/GRE/CG/",F5),Sheet1!$A$1:$B$1250,2,FALSE)'
$wsHoja1.Range("H5").Value = 1

# Row 6
$wsHoja1.Range("A6").Value = 4
$wsHoja1.Range("B6").Value = 143
$wsHoja1.Range("C6").Value = 45061
$wsHoja1.Range("D6").Value = "Seguimiento"
$wsHoja1.Range("E6").Value = "En el SAI"
$wsHoja1.Range("F6").Formula = '=CONCATENATE(TEXT(B6,"0000"),"/",TEXT(H6,"00"))'
$wsHoja1.Range("G6").Formula = '=VLOOKUP(CONCATENATE("PU/GRE/CG/",F6),Sheet1!$A$1:$B$1250,2,FALSE)'
$wsHoja1.Range("H6").Value = 1

# Row 7
$wsHoja1.Range("A7").Value = 5
$wsHoja1.Range("B7").Value = 256
$wsHoja1.Range("C7").Value = 45058
$wsHoja1.Range("D7").Value = "Seguimiento "
$wsHoja1.Range("E7").Value = "En el SAI"
$wsHoja1.Range("F7").Formula = '=CONCATENATE(TEXT(B7,"0000"),"/",TEXT(H7,"00"))'
$wsHoja1.Range("G7").Formula = '=VLOOKUP(CONCATENATE("PU/GRE/CG/",F7),Sheet1!$A$1:$B$1250,2,FALSE)'
$wsHoja1.Range("H7").Value = 1

# ---------------------------------------------------------------------------
# 3) View state: the author ended up scrolled near the bottom of the lookup
#    sheet (having just pasted the new ids) before returning focus to the
#    main sheet with E1 selected.
# ---------------------------------------------------------------------------

$wsLookup.Activate()
$wsLookup.Range("B1136").Select()

$wsHoja1.Activate()
$wsHoja1.Range("E1").Select()
